# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.506.37'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '2.379.83'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '''314.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").Value = '''108.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.59%  '
$ws.Range("D7").Value = '''0.631'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.39%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").Value = '''0.614'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.14%  '
$ws.Range("D10").Value = '''41.01'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.01%  '
$ws.Range("D11").Value = '''0.0922'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("D12").Value = '''8.56'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("E13").Value = '  +0.91%  '
$ws.Range("D14").Value = '''0.987'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("D15").Value = '2.739.44'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").Value = '''15.39'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.36%  '
$ws.Range("D17").Value = '2.362.66'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = '45.482.22'
$ws.Range("E18").Value = '  +1.21%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '''7.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.30%  '
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").Value = '''13.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.52%  '
$ws.Range("D21").Value = '''0.0000107'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").Value = '''73.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.03%  '
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("D24").Value = '''260.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.20%  '
$ws.Range("D25").Value = '''2.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("D27").Value = '''11.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  -5.44%  '
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("D30").Value = '''0.0985'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.06%  '
$ws.Range("D31").Value = '''22.43'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.09%  '
$ws.Range("D32").Value = '''37.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.41%  '
$ws.Range("D33").Value = '''166.63'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("D34").Value = '''2.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.92%  '
$ws.Range("E35").Value = '  -1.34%  '
$ws.Range("D36").Value = '''0.119'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("D37").Value = '''4.71'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.15%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = '''1.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.57%  '
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = '''4.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.25%  '
$ws.Range("D40").Value = '''2.99'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").Value = '''0.0358'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("D42").Value = '''99.03'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.63%  '
$ws.Range("D43").Value = '''69.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.24%  '
$ws.Range("D44").Value = '''0.230'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.43%  '
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = '''12.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.67%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Value = '''1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").Value = '1.822.57'
$ws.Range("E47").Value = '  +9.61%  '
$ws.Range("D48").Value = '''84.24'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.53%  '
$ws.Range("D49").Value = '''5.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.49%  '
$ws.Range("E50").Value = '  +2.69%  '
$ws.Range("D51").Value = '''111.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.69%  '
